$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was added to the top of the data set (row 30).
# This shifts every existing data row (30-139) down by one row, so insert
# a fresh row at 30 first and then populate it with the new record - a
# near-duplicate of the (now shifted-down) former row 30, but with an
# updated date (Fecha) and volume (Volumen).
$ws.Rows(30).Insert()

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C30").Value = "Metropolitana"
$ws.Range("D30").Value = 45012
$ws.Range("E30").Value = 13
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100104
$ws.Range("H30").Value = "Frutos de pepita"
$ws.Range("I30").Value = 100104003
$ws.Range("J30").Value = "Membrillo"
$ws.Range("K30").Value = "Champion"
$ws.Range("L30").Value = "Especial"
$ws.Range("M30").Value = 6
$ws.Range("N30").Value = 280000
$ws.Range("O30").Value = 280000
$ws.Range("P30").Value = 280000
$ws.Range("Q30").Value = "$/bins (450 kilos)"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 622
$ws.Range("T30").Value = 450
